$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.490.25'
$ws.Range('E2').Value = '  -0.66%  '
$ws.Range('D3').Value = '1.621.33'
$ws.Range('E3').Value = '  -1.13%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = "'211.50"
$ws.Range('E5').Value = '  -0.65%  '
$ws.Range('E6').Value = '  -0.83%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = "'23.09"
$ws.Range('E8').Value = '  -0.23%  '
$ws.Range('E9').Value = '  +1.58%  '
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('D11').Value = "'0.0881"
$ws.Range('E11').Value = '  -1.61%  '
$ws.Range('D12').Value = '1.850.84'
$ws.Range('E12').Value = '  -1.13%  '
$ws.Range('D13').Value = '1.641.74'
$ws.Range('E13').Value = '  +0.03%  '
$ws.Range('E14').Value = '  -0.26%  '
$ws.Range('E15').Value = '  -2.04%  '
$ws.Range('D16').Value = "'65.33"
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('D17').Value = '27.480.55'
$ws.Range('E17').Value = '  -0.61%  '
$ws.Range('D18').Value = "'229.21"
$ws.Range('E18').Value = '  -0.56%  '
$ws.Range('D19').Value = '0.0₃0718'
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('E20').Value = '  -2.57%  '
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('D22').Value = "'10.45"
$ws.Range('E22').Value = '  +4.11%  '
$ws.Range('E23').Value = '  +1.16%  '
$ws.Range('E24').Value = '  +8.15%  '
$ws.Range('D25').Value = "'149.26"
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('D26').Value = "'6.87"
$ws.Range('E26').Value = '  -0.98%  '
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('D29').Value = "'15.50"
$ws.Range('E29').Value = '  -0.90%  '
$ws.Range('E30').Value = '  -0.64%  '
$ws.Range('E31').Value = '  -0.64%  '
$ws.Range('E32').Value = '  -1.01%  '
$ws.Range('D33').Value = '1.465.13'
$ws.Range('E33').Value = '  +1.43%  '
$ws.Range('E34').Value = '  -2.38%  '
$ws.Range('E35').Value = '  -1.73%  '
$ws.Range('E36').Value = '  -1.85%  '
$ws.Range('D37').Value = "'0.948"
$ws.Range('E37').Value = '  +5.24%  '
$ws.Range('E38').Value = '  +0.18%  '
$ws.Range('E39').Value = '  -0.59%  '
$ws.Range('D40').Value = "'0.552"
$ws.Range('E40').Value = '  -2.64%  '
$ws.Range('B41').Value = 'WEMIXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').Value = "'1.02"
$ws.Range('E41').Value = '  -0.84%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = "'1.00"
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').Value = "'67.94"
$ws.Range('E43').Value = '  -3.36%  '
$ws.Range('E44').Value = '  +0.79%  '
$ws.Range('E45').Value = '  -2.03%  '
$ws.Range('D46').Value = "'5.33"
$ws.Range('E46').Value = '  -4.98%  '
$ws.Range('E47').Value = '  +2.27%  '
$ws.Range('D48').Value = '1.761.24'
$ws.Range('E48').Value = '  -1.18%  '
$ws.Range('D49').Value = "'87.09"
$ws.Range('E49').Value = '  +1.15%  '
$ws.Range('E50').Value = '  -0.76%  '
$ws.Range('E51').Value = '  +0.37%  '

# Reset style on cells that needed the text quote-prefix, to avoid a stray quotePrefix style
$ws.Range('D5').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
